# Add a new "2023" column (S) to the table, mirroring the existing
# formatting of the last data column (R) for the header row and the two
# data rows, then fill in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column R's formatting (rows 4-6, the only rows with data in R)
# into column S so the new column looks consistent with the rest of the
# table (borders, number format, font, etc.).
$ws.Range("R4:R6").Copy()
$ws.Range("S4:S6").PasteSpecial(-4122)

# Write the new year header and the two new data points.
$ws.Range("S4").Value = 2023
$ws.Range("S5").Value = 7.1262361838278068
$ws.Range("S6").Value = 10.974456007568591
